$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (shifts existing rows 10-25 down to 11-26)
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with the new weekly data point
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = "9/8/2023"
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 100112001
$ws.Range("G10").Value = "Berenjena"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 160
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 5500
$ws.Range("M10").Value = 5250
$ws.Range("N10").Value = "$/caja 60 unidades"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 88
$ws.Range("Q10").Value = 60
$ws.Range("R10").Value = "Hortaliza"
